# Generate Report for Handoff
# Replaces the source-file GUID token (8d142fa8-...) with a new one
# (05aa3b24-...), refreshes the handoff/handback timestamps, and clears
# out the "Latest Target File" / "Latest Handback File" columns (the
# handback for the new run hasn't happened yet) on both locale sheets.

$wb = $excel.ActiveWorkbook

$newGuid         = "05aa3b24-e7d8-42c2-bc5d-2d0d1660238f"
$oldGuidFile     = "8d142fa8-10d9-420a-8baa-1aa05f6bca52.md"
$newGuidFile     = "05aa3b24-e7d8-42c2-bc5d-2d0d1660238f.md"
$oldGuidPath     = "e2e\8d142fa8-10d9-420a-8baa-1aa05f6bca52.md"
$newGuidPath     = "e2e\05aa3b24-e7d8-42c2-bc5d-2d0d1660238f.md"

$oldHash         = "72edb2ab0556393848e9ad449f04c7741384edba"
$newHash         = "0b917953ecef39e11a548e8e49c376282f6b3739"

$sourceHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/oltest/blob/4b7e7cfd705e507c8461dcfa06335b7d3e385070/e2e/8d142fa8-10d9-420a-8baa-1aa05f6bca52.md"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = $newGuidFile
$ov.Range("B2").Value = $newGuidPath
$ov.Range("G2").Value = "2016-08-12 23:16:24"

$ovAddr = $sourceHyperlinkAddress
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), $ovAddr, "", "", $newGuidPath)

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = $newGuidFile
$zh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-12 23:16:16"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Range("I2").Style = "Normal"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $sourceHyperlinkAddress, "", "", $newGuidFile)

$zh.Columns.Item(9).ColumnWidth = 17.8
$zh.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = $newGuidFile
$de.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$de.Range("H2").Value = "2016-08-12 23:16:24"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Range("I2").Style = "Normal"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $sourceHyperlinkAddress, "", "", $newGuidFile)

$de.Columns.Item(9).ColumnWidth = 17.8
$de.Columns.Item(10).ColumnWidth = 20.8
